$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelas ODS")
$ws.Activate()
$win = $excel.ActiveWindow
Write-Host ("TopLeftCell? trying property")
try {
  $ws.Range("A540").Select()
  $win.ScrollRow = 540
  $win.ScrollColumn = 1
} catch {
  Write-Host ("Error: " + $_)
}
$ws.Range("C560").Select()
Write-Host "done"
